$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New camp row (row 5) appended beneath the existing 3 camp rows.
$ws.Cells.Item(5, 1).Value = "tester2"
$ws.Cells.Item(5, 2).Value = "01-01-01"
$ws.Cells.Item(5, 3).Value = "02-02-02"
$ws.Cells.Item(5, 4).Value = "SCSE"
$ws.Cells.Item(5, 5).Value = "ntu"
$ws.Cells.Item(5, 6).Value = 100.0
$ws.Cells.Item(5, 7).Value = 10.0
$ws.Cells.Item(5, 8).Value = "tester2"
$ws.Cells.Item(5, 9).Value = $true
$ws.Cells.Item(5, 10).Value = "HUKUMAR"
$ws.Cells.Item(5, 11).Value = ""
$ws.Cells.Item(5, 12).Value = ""

# New trailing column (M / 13) used for an extra (empty) field on every
# data row, including the newly appended one.
$ws.Cells.Item(2, 13).Value = ""
$ws.Cells.Item(3, 13).Value = ""
$ws.Cells.Item(4, 13).Value = ""
$ws.Cells.Item(5, 13).Value = ""
